# Fix the two "usercontent" field markers so they use the "long" Word
# field-code syntax (a run with <w:fldChar begin>, a run holding the
# instruction text in <w:instrText>, a run with <w:fldChar separate>
# and a run with <w:fldChar end>) instead of the compact
# <w:fldSimple w:instr="..."/> form. Word normally performs exactly
# this expansion as soon as such a field is touched/edited and
# re-saved; we reproduce that normalization here.

$d = $word.ActiveDocument

function Expand-Instr($rawInstr) {
    # Un-escape the XML entities WordOpenXML gives us back, then
    # re-escape for safe embedding inside the new <w:instrText>.
    $unescaped = $rawInstr.Replace("&quot;", '"').Replace("&apos;", "'").Replace("&lt;", "<").Replace("&gt;", ">").Replace("&amp;", "&")
    return $unescaped.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

$fldPattern = '<w:fldSimple w:instr="[^"]*"\s*/>|<w:fldSimple w:instr="[^"]*"\s*>\s*</w:fldSimple>'
$instrPattern = '<w:fldSimple w:instr="([^"]*)"'

# Process paragraphs from last to first so replacing one doesn't shift
# the character offsets of the ones we still have to visit.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range

    $openXml = $full.WordOpenXML
    if ($openXml -notmatch '<w:fldSimple') { continue }

    # WordOpenXML wraps the requested range together with trailing
    # document scaffolding; pull out just this paragraph's markup.
    if ($openXml -notmatch '(?s)<w:body>(<w:p\b.*?</w:p>)') { continue }
    $paraXml = $matches[1]

    if ($paraXml -notmatch $fldPattern) { continue }

    # A paragraph could in principle hold more than one fldSimple;
    # expand every one of them, left to right.
    $newParaXml = $paraXml
    while ($newParaXml -match $fldPattern) {
        $fldXml = $matches[0]
        if ($fldXml -notmatch $instrPattern) { break }
        $instr = Expand-Instr $matches[1]

        $idx = $newParaXml.IndexOf($fldXml)
        $prefix = $newParaXml.Substring(0, $idx)
        $suffix = $newParaXml.Substring($idx + $fldXml.Length)

        $expanded = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' + `
            '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' + `
            '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' + `
            '<w:r><w:fldChar w:fldCharType="end"/></w:r>'

        $newParaXml = $prefix + $expanded + $suffix
    }

    $xmlPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    # Replace everything up to (but not including) the paragraph mark
    # so the paragraph's own properties stay attached to it.
    $target = $d.Range($full.Start, $full.End - 1)
    $target.InsertXML($xmlPkg)
}
